# Absenzenlisten-Template 2016/2017 - fix oversized font in a few
# table cells ("zu grosse Schrift in einzelnen Zellen korrigiert").
#
# The attendance table has one header block followed by 23 student
# rows. In each student row, three particular cells (two empty
# "absence" cells and the first "X" marker cell of the following
# month block) were missing the explicit 10pt run formatting that
# all the sibling cells already carry, so Word rendered them at the
# default 11pt - visibly larger than the rest of the row. We restore
# the missing <w:sz w:val="20"/><w:szCs w:val="20"/> (10pt) run
# properties on the paragraph mark (and on the run itself for the
# "X" cell) to match the rest of the table.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Student rows are table rows 3..25 (row 1 = month header, row 2 =
# column header). Within each row, cells 6 and 7 are the two empty
# cells, and cell 18 holds the "X" run - all three lack the 10pt
# formatting present on every other cell in the row.
$targetCells = 6, 7, 18

for ($rowIndex = 3; $rowIndex -le 25; $rowIndex++) {
    $row = $table.Rows.Item($rowIndex)
    foreach ($cellIndex in $targetCells) {
        $cell = $row.Cells.Item($cellIndex)
        $cell.Range.Font.Size = 10
        $cell.Range.Font.SizeBi = 10
    }
}

Write-Output "Fixed font size on $($targetCells.Count) cells across 23 student rows."
